$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I19").Value = -0.5348266449214059
$ws.Range("J19").Value = 0.1282018702557457
$ws.Range("K19").Value = -0.2262958989451835
$ws.Range("L19").Value = 2.013425631032702
